$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 40848.5
$ws.Range("J3").Value = 40848.5
$ws.Range("L3").Value = 40848.5
$ws.Range("N3").Value = -41076.5

$ws.Range("H70").Value = 1503.2858
$ws.Range("I70").Value = 1088
$ws.Range("J70").Value = 1734
$ws.Range("K70").Value = 3264
$ws.Range("L70").Value = 5202
$ws.Range("M70").Value = -2994
$ws.Range("N70").Value = -5742

$ws.Range("H73").Value = 1503.2858
$ws.Range("I73").Value = 1088
$ws.Range("J73").Value = 1734
$ws.Range("K73").Value = 3264
$ws.Range("L73").Value = 5202
$ws.Range("M73").Value = -2328
$ws.Range("N73").Value = -7074

$ws.Range("H102").Value = 40848.5
$ws.Range("J102").Value = 40848.5
$ws.Range("L102").Value = 40848.5
$ws.Range("N102").Value = -47338.5

$ws.Range("H138").Value = 4312498.5
$ws.Range("J138").Value = 15629864
$ws.Range("L138").Value = 46889592
$ws.Range("N138").Value = -46899872

$ws.Range("H141").Value = 1764.9524
$ws.Range("I141").Value = 1671.5264
$ws.Range("J141").Value = 2652.5
$ws.Range("K141").Value = 5014.5792
$ws.Range("L141").Value = 7957.5
$ws.Range("M141").Value = 165.4207999999999
$ws.Range("N141").Value = -18317.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I74").Value = 12821287
$ws.Range("J74").Value = 7138.6665
$ws.Range("K74").Value = 12821287
$ws.Range("L74").Value = 7138.6665
$ws.Range("M74").Value = -12820413
$ws.Range("N74").Value = -8886.666499999999

$ws.Range("I77").Value = 12821287
$ws.Range("J77").Value = 7138.6665
$ws.Range("K77").Value = 64106435
$ws.Range("L77").Value = 35693.3325
$ws.Range("M77").Value = -64102067
$ws.Range("N77").Value = -44429.3325

$ws.Range("H102").Value = 4691.5806
$ws.Range("I102").Value = 5349.522
$ws.Range("J102").Value = 2800
$ws.Range("K102").Value = 5349.522
$ws.Range("L102").Value = 2800
$ws.Range("M102").Value = -3727.522
$ws.Range("N102").Value = -6044

$ws.Range("H122").Value = 38004
$ws.Range("I122").Value = 55506
$ws.Range("K122").Value = 166518
$ws.Range("M122").Value = -164068

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1749.6666
$ws.Range("I58").Value = 1070.2222
$ws.Range("J58").Value = 3278.4167
$ws.Range("K58").Value = 1070.2222
$ws.Range("L58").Value = 3278.4167
$ws.Range("M58").Value = -867.2221999999999
$ws.Range("N58").Value = -3684.4167

$ws.Range("H136").Value = 1749.6666
$ws.Range("I136").Value = 1070.2222
$ws.Range("J136").Value = 3278.4167
$ws.Range("K136").Value = 3210.6666
$ws.Range("L136").Value = 9835.250100000001
$ws.Range("M136").Value = -660.6665999999996
$ws.Range("N136").Value = -14935.2501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 5636.4443
$ws.Range("I62").Value = 2333.3333
$ws.Range("K62").Value = 6999.999899999999
$ws.Range("M62").Value = -6313.999899999999

$ws.Range("H63").Value = 3854.1667
$ws.Range("I63").Value = 737
$ws.Range("J63").Value = 6971.3335
$ws.Range("K63").Value = 2211
$ws.Range("L63").Value = 20914.0005
$ws.Range("M63").Value = -1462
$ws.Range("N63").Value = -22412.0005

$ws.Range("H64").Value = 2842.25
$ws.Range("I64").Value = 941.3333
$ws.Range("J64").Value = 3982.8
$ws.Range("K64").Value = 2823.9999
$ws.Range("L64").Value = 11948.4
$ws.Range("M64").Value = -2553.9999
$ws.Range("N64").Value = -12488.4

$ws.Range("H65").Value = 5636.4443
$ws.Range("I65").Value = 2333.3333
$ws.Range("K65").Value = 20999.9997
$ws.Range("M65").Value = -17567.9997

$ws.Range("H66").Value = 3854.1667
$ws.Range("I66").Value = 737
$ws.Range("J66").Value = 6971.3335
$ws.Range("K66").Value = 6633
$ws.Range("L66").Value = 62742.0015
$ws.Range("M66").Value = -2889
$ws.Range("N66").Value = -70230.0015

$ws.Range("H67").Value = 2842.25
$ws.Range("I67").Value = 941.3333
$ws.Range("J67").Value = 3982.8
$ws.Range("K67").Value = 2823.9999
$ws.Range("L67").Value = 11948.4
$ws.Range("M67").Value = -1887.9999
$ws.Range("N67").Value = -13820.4

$ws.Range("H68").Value = 1067.661
$ws.Range("J68").Value = 1270.4117
$ws.Range("L68").Value = 3811.2351
$ws.Range("N68").Value = -5433.2351

$ws.Range("H71").Value = 1067.661
$ws.Range("J71").Value = 1270.4117
$ws.Range("L71").Value = 11433.7053
$ws.Range("N71").Value = -19545.7053

$ws.Range("H107").Value = 922.07275
$ws.Range("I107").Value = 316.2963
$ws.Range("J107").Value = 1506.2142
$ws.Range("K107").Value = 948.8888999999999
$ws.Range("L107").Value = 4518.642599999999
$ws.Range("M107").Value = 971.1111000000001
$ws.Range("N107").Value = -8358.642599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 8849.833000000001
$ws.Range("I102").Value = 12149.75
$ws.Range("J102").Value = 2250
$ws.Range("K102").Value = 12149.75
$ws.Range("L102").Value = 2250
$ws.Range("M102").Value = -10527.75
$ws.Range("N102").Value = -5494

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 806.2381
$ws.Range("I22").Value = 797.5
$ws.Range("J22").Value = 811.61536
$ws.Range("K22").Value = 797.5
$ws.Range("L22").Value = 811.61536
$ws.Range("M22").Value = -502.5
$ws.Range("N22").Value = -1401.61536

$ws.Range("H27").Value = 806.2381
$ws.Range("I27").Value = 797.5
$ws.Range("J27").Value = 811.61536
$ws.Range("K27").Value = 797.5
$ws.Range("L27").Value = 811.61536
$ws.Range("M27").Value = -690.5
$ws.Range("N27").Value = -1025.61536

$ws.Range("H55").Value = 624.5
$ws.Range("I55").Value = 316.66666
$ws.Range("J55").Value = 727.1111
$ws.Range("K55").Value = 316.66666
$ws.Range("L55").Value = 727.1111
$ws.Range("M55").Value = -143.66666
$ws.Range("N55").Value = -1073.1111

$ws.Range("H94").Value = 46871.25
$ws.Range("J94").Value = 46871.25
$ws.Range("L94").Value = 46871.25
$ws.Range("N94").Value = -48223.25

$ws.Range("H122").Value = 7928.9414
$ws.Range("I122").Value = 10857.714
$ws.Range("K122").Value = 32573.142
$ws.Range("M122").Value = -30123.142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 10000
$ws.Range("J22").Value = 10000
$ws.Range("L22").Value = 10000
$ws.Range("N22").Value = -10586

$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344

$ws.Range("H136").Value = 2778654
$ws.Range("I136").Value = 3031119.5
$ws.Range("J136").Value = 1533.3334
$ws.Range("K136").Value = 9093358.5
$ws.Range("L136").Value = 4600.0002
$ws.Range("M136").Value = -9090808.5
$ws.Range("N136").Value = -9700.0002
